$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the listed rows
$ws.Range("F3").Value = -2
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = -2
$ws.Range("F13").Value = 1
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = -13
$ws.Range("F17").Value = -6
$ws.Range("F19").Value = -5
$ws.Range("F20").Value = 9
$ws.Range("F21").Value = -3
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 4
